$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 6055
$ws.Range("I43").Value = 1244.5
$ws.Range("J43").Value = 7979.2
$ws.Range("K43").Value = 1244.5
$ws.Range("L43").Value = 7979.2
$ws.Range("M43").Value = -1175.5
$ws.Range("N43").Value = -8117.2
# Row 53
$ws.Range("H53").Value = 10889.263
$ws.Range("I53").Value = 383.63635
$ws.Range("J53").Value = 25334.5
$ws.Range("K53").Value = 383.63635
$ws.Range("L53").Value = 25334.5
$ws.Range("M53").Value = 253.36365
$ws.Range("N53").Value = -26608.5
# Row 76
$ws.Range("H76").Value = 4875.6665
$ws.Range("I76").Value = 3702.4167
$ws.Range("J76").Value = 6440
$ws.Range("K76").Value = 3702.4167
$ws.Range("L76").Value = 6440
$ws.Range("M76").Value = -3387.4167
$ws.Range("N76").Value = -7070
# Row 79
$ws.Range("H79").Value = 4875.6665
$ws.Range("I79").Value = 3702.4167
$ws.Range("J79").Value = 6440
$ws.Range("K79").Value = 3702.4167
$ws.Range("L79").Value = 6440
$ws.Range("M79").Value = -2610.4167
$ws.Range("N79").Value = -8624
# Row 92
$ws.Range("H92").Value = 1061.0714
$ws.Range("I92").Value = 331.65
$ws.Range("J92").Value = 2884.625
$ws.Range("K92").Value = 331.65
$ws.Range("L92").Value = 2884.625
$ws.Range("M92").Value = 916.35
$ws.Range("N92").Value = -5380.625
# Row 134
$ws.Range("H134").Value = 193665.22
$ws.Range("J134").Value = 193665.22
$ws.Range("L134").Value = 193665.22
$ws.Range("N134").Value = -203805.22

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9571.223
$ws.Range("I32").Value = 6333.3403
$ws.Range("K32").Value = 6333.3403
$ws.Range("M32").Value = -6046.3403
# Row 74
$ws.Range("H74").Value = 108816.78
$ws.Range("I74").Value = 85775.336
$ws.Range("J74").Value = 154899.67
$ws.Range("K74").Value = 85775.336
$ws.Range("L74").Value = 154899.67
$ws.Range("M74").Value = -84901.336
$ws.Range("N74").Value = -156647.67
# Row 77
$ws.Range("H77").Value = 108816.78
$ws.Range("I77").Value = 85775.336
$ws.Range("J77").Value = 154899.67
$ws.Range("K77").Value = 428876.68
$ws.Range("L77").Value = 774498.3500000001
$ws.Range("M77").Value = -424508.68
$ws.Range("N77").Value = -783234.3500000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2782919.8
$ws.Range("I86").Value = 4002117.5
$ws.Range("J86").Value = 12015.363
$ws.Range("K86").Value = 4002117.5
$ws.Range("L86").Value = 12015.363
$ws.Range("M86").Value = -4000994.5
$ws.Range("N86").Value = -14261.363
# Row 89
$ws.Range("H89").Value = 2782919.8
$ws.Range("I89").Value = 4002117.5
$ws.Range("J89").Value = 12015.363
$ws.Range("K89").Value = 20010587.5
$ws.Range("L89").Value = 60076.815
$ws.Range("M89").Value = -20004971.5
$ws.Range("N89").Value = -71308.815
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 1500
$ws.Range("I23").Value = 1500
$ws.Range("K23").Value = 1500
$ws.Range("M23").Value = -1260
# Row 27
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 1500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1308
# Row 28
$ws.Range("H28").Value = 28010.75
$ws.Range("J28").Value = 28010.75
$ws.Range("L28").Value = 28010.75
$ws.Range("N28").Value = -28500.75
# Row 69
$ws.Range("H69").Value = 38636
$ws.Range("I69").Value = 23295
$ws.Range("K69").Value = 23295
$ws.Range("M69").Value = -22546
# Row 72
$ws.Range("H72").Value = 38636
$ws.Range("I72").Value = 23295
$ws.Range("K72").Value = 69885
$ws.Range("M72").Value = -66141
# Row 86
$ws.Range("H86").Value = 7836.654
$ws.Range("I86").Value = 6605.8237
$ws.Range("K86").Value = 6605.8237
$ws.Range("M86").Value = -5482.8237
# Row 89
$ws.Range("H89").Value = 7836.654
$ws.Range("I89").Value = 6605.8237
$ws.Range("K89").Value = 33029.1185
$ws.Range("M89").Value = -27413.1185
# Row 134
$ws.Range("H134").Value = 24973.834
$ws.Range("I134").Value = 28187.838
$ws.Range("K134").Value = 84563.514
$ws.Range("M134").Value = -82028.514
# Row 135
$ws.Range("H135").Value = 109822.8
$ws.Range("J135").Value = 109822.8
$ws.Range("L135").Value = 109822.8
$ws.Range("N135").Value = -119962.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 1689.75
$ws.Range("I3").Value = 1179.6666
$ws.Range("K3").Value = 3538.9998
$ws.Range("M3").Value = -3426.9998
# Row 117
$ws.Range("H117").Value = 1887.2222
$ws.Range("I117").Value = 1160
$ws.Range("J117").Value = 2095
$ws.Range("K117").Value = 3480
$ws.Range("L117").Value = 6285
$ws.Range("M117").Value = -38
$ws.Range("N117").Value = -13169
# Row 129
$ws.Range("H129").Value = 1334182.9
$ws.Range("I129").Value = 1818661.9
$ws.Range("J129").Value = 1865.5
$ws.Range("K129").Value = 5455985.699999999
$ws.Range("L129").Value = 5596.5
$ws.Range("M129").Value = -5450985.699999999
$ws.Range("N129").Value = -15596.5
# Row 130
$ws.Range("H130").Value = 2930.6667
$ws.Range("I130").Value = 2710.8572
$ws.Range("J130").Value = 3700
$ws.Range("K130").Value = 8132.571599999999
$ws.Range("L130").Value = 11100
$ws.Range("M130").Value = -3112.571599999999
$ws.Range("N130").Value = -21140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 23000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
# Row 70
$ws.Range("H70").Value = 33336872
$ws.Range("I70").Value = 50003252
$ws.Range("K70").Value = 50003252
$ws.Range("M70").Value = -50002982
# Row 73
$ws.Range("H73").Value = 33336872
$ws.Range("I73").Value = 50003252
$ws.Range("K73").Value = 50003252
$ws.Range("M73").Value = -50002316
# Row 133
$ws.Range("H133").Value = 84390
$ws.Range("J133").Value = 84390
$ws.Range("L133").Value = 84390
$ws.Range("N133").Value = -94510

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 106
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524
# Row 132
$ws.Range("H132").Value = 15498.258
$ws.Range("I132").Value = 17953.875
$ws.Range("J132").Value = 7079
$ws.Range("K132").Value = 53861.625
$ws.Range("L132").Value = 21237
$ws.Range("M132").Value = -51331.625
$ws.Range("N132").Value = -26297
# Row 136
$ws.Range("H136").Value = 53630.465
$ws.Range("I136").Value = 120950
$ws.Range("J136").Value = 5945.7915
$ws.Range("K136").Value = 362850
$ws.Range("L136").Value = 17837.3745
$ws.Range("M136").Value = -360300
$ws.Range("N136").Value = -22937.3745
# Row 139
$ws.Range("H139").Value = 68774.5
$ws.Range("J139").Value = 68742.28999999999
$ws.Range("L139").Value = 68742.28999999999
$ws.Range("N139").Value = -79022.28999999999
# Row 140
$ws.Range("H140").Value = 83899.60000000001
$ws.Range("J140").Value = 97374.5
$ws.Range("L140").Value = 97374.5
$ws.Range("N140").Value = -107734.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 7250501.5
$ws.Range("I81").Value = 13891979
$ws.Range("J81").Value = 5253
$ws.Range("K81").Value = 27783958
$ws.Range("L81").Value = 10506
$ws.Range("M81").Value = -27782897
$ws.Range("N81").Value = -12628
# Row 84
$ws.Range("H84").Value = 7250501.5
$ws.Range("I84").Value = 13891979
$ws.Range("J84").Value = 5253
$ws.Range("K84").Value = 138919790
$ws.Range("L84").Value = 52530
$ws.Range("M84").Value = -138914486
$ws.Range("N84").Value = -63138
# Row 107
$ws.Range("H107").Value = 37038290
$ws.Range("I107").Value = 90909790
$ws.Range("K107").Value = 272729370
$ws.Range("M107").Value = -272727450
# Row 132
$ws.Range("H132").Value = 24079748
$ws.Range("I132").Value = 28581576
$ws.Range("K132").Value = 85744728
$ws.Range("M132").Value = -85744728
